# The sheet tracks one weekly price observation per row (rows 2..250).
# This edit adds a new, more recent observation ("Fruta / hortaliza, semanal")
# by inserting a fresh row at position 154 and shifting every row from the
# old 154 through 250 down by one (they land on 155..251 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 154; rows 154..250 shift down to 155..251.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A154").Value = 4
$ws.Range("B154").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C154").Value = "Los Lagos"
$ws.Range("D154").Value = 44603
$ws.Range("E154").Value = 10
$ws.Range("F154").Value = 100114014
$ws.Range("G154").Value = "Betarraga"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 1000
$ws.Range("K154").Value = 900
$ws.Range("L154").Value = 1000
$ws.Range("M154").Value = 950
$ws.Range("N154").Value = "$/paquete 5 unidades"
$ws.Range("O154").Value = "Región del Maule"
$ws.Range("P154").Value = 190
$ws.Range("Q154").Value = 5
$ws.Range("R154").Value = "Hortaliza"
